$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dates/quality/volume/prices updated
$ws.Range("D2").Value = Get-Date -Year 2021 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2364
$ws.Range("N2").Value = "`$/docena de matas"
$ws.Range("P2").Value = 394
$ws.Range("Q2").Value = 6

# Row 3: date/volume/prices updated
$ws.Range("D3").Value = Get-Date -Year 2021 -Month 6 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = 1080
$ws.Range("P3").Value = 180

# Row 5: date/quality/volume/prices updated
$ws.Range("D5").Value = Get-Date -Year 2022 -Month 3 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("N5").Value = "`$/paquete"
$ws.Range("P5").Value = 1900
$ws.Range("Q5").Value = 1
